# Thêm khách hàng mới "Bạch Nhi" vào đầu danh sách (dòng 2), đẩy các dòng
# hiện có xuống 1 dòng.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Chèn một dòng trống mới tại vị trí dòng 2 (các dòng bên dưới dịch xuống).
$ws.Rows.Item(2).Insert()

# Điền dữ liệu khách hàng mới vào dòng 2 vừa chèn.
$ws.Range("A2").Value = "KH"
$ws.Range("B2").Value = 409
$ws.Range("C2").Value = "Bạch Nhi"
$ws.Range("D2").Value = "CẦN THƠ"
$ws.Range("I2").Value = 5000000
$ws.Range("J2").Value = 10000000
